$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.514.09"
$ws.Range("E2").Value = "  +0.50%  "

$ws.Range("D3").Value = "3.764.81"
$ws.Range("E3").Value = "  -0.66%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.32%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "621.52"
$ws.Range("E5").Value = "  +0.74%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "181.77"
$ws.Range("E6").Value = "  +2.59%  "

$ws.Range("D7").Value = "3.759.11"
$ws.Range("E7").Value = "  -0.66%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.21%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.536"
$ws.Range("E9").Value = "  -2.25%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.169"
$ws.Range("E10").Value = "  +0.77%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.40"
$ws.Range("E11").Value = "  -0.08%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.485"
$ws.Range("E12").Value = "  -3.32%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "40.65"
$ws.Range("E13").Value = "  +0.15%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000258"
$ws.Range("E14").Value = "  +0.00%  "

$ws.Range("D15").Value = "4.375.36"
$ws.Range("E15").Value = "  -1.35%  "

$ws.Range("D16").Value = "3.747.25"
$ws.Range("E16").Value = "  -1.49%  "

$ws.Range("D17").Value = "70.325.61"
$ws.Range("E17").Value = "  -0.04%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.64"
$ws.Range("E18").Value = "  +0.38%  "

$ws.Range("E19").Value = "  -1.98%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "508.95"
$ws.Range("E20").Value = "  -3.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "16.58"
$ws.Range("E21").Value = "  -0.80%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.27"
$ws.Range("E22").Value = "  -2.39%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.727"
$ws.Range("E23").Value = "  -2.47%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.58"
$ws.Range("E24").Value = "  +3.37%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "87.35"
$ws.Range("E25").Value = "  -1.42%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.09"
$ws.Range("E26").Value = "  -3.37%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.26"
$ws.Range("E27").Value = "  +2.76%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0000136"
$ws.Range("E28").Value = "  +10.36%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.50"
$ws.Range("E30").Value = "  -0.76%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.94"
$ws.Range("E31").Value = "  +2.30%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.91"
$ws.Range("E32").Value = "  -0.19%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.77"
$ws.Range("E33").Value = "  -4.30%  "

$ws.Range("E34").Value = "  -0.12%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.19%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.05"
$ws.Range("E36").Value = "  +0.35%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.15"
$ws.Range("E37").Value = "  -1.27%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.350"
$ws.Range("E38").Value = "  +1.92%  "

$ws.Range("E39").Value = "  +4.61%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.24"
$ws.Range("E40").Value = "  +17.56%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.10"
$ws.Range("E41").Value = "  -3.13%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "49.98"
$ws.Range("E42").Value = "  -3.32%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "432.48"
$ws.Range("E43").Value = "  +1.05%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "44.57"
$ws.Range("E44").Value = "  +0.08%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.67"
$ws.Range("E45").Value = "  -2.24%  "

$ws.Range("D46").Value = "2.972.05"
$ws.Range("E46").Value = "  -5.54%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0364"
$ws.Range("E47").Value = "  -1.21%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.43"
$ws.Range("E48").Value = "  -1.48%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "136.62"
$ws.Range("E50").Value = "  -0.20%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.49"
$ws.Range("E51").Value = "  -1.71%  "
